# VerveStacks_DEU_grids model update - 2025-08-19
# Add a new "Grid" process set (IRE / g[_]*) row to the VEDA_Sets-Proc sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# New row 21: PSET_SET=IRE, PSET_PN=g[_]*, SetName=Grid
$ws.Range("F21").Value = "Grid"
$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"

# Leave the cursor on the newly entered cell, as in the source edit.
$ws.Range("B21").Select() | Out-Null
